$d = $word.ActiveDocument

# Locate the paragraph that holds the " {Name}" placeholder (the one that
# currently carries the "_GoBack" bookmark at its start) and insert a new
# paragraph right after it for the new "{Total}" placeholder.
$namePara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text -match "\{Name\}") {
        $namePara = $p
    }
}

$namePara.Range.InsertParagraphAfter()

$totalPara = $namePara.Next()
$totalPara.Range.Text = "{Total}"

# Move the "_GoBack" bookmark so it now sits inside the new paragraph,
# between "{Total" and "}" (mirrors where it used to sit, just before the
# closing brace, in the original " {Name}" paragraph).
$insertPos = $totalPara.Range.Start + 6
$bmRange = $d.Range($insertPos, $insertPos)
$d.Bookmarks.Add("_GoBack", $bmRange)
